# Update the "Levels" row (row 4) of the RQ134 PH2 variable-definitions table:
#  - Datasource levels now list the individual Phase-2 data collection waves
#  - Attribute levels gain the new "Personal Safety / Affiliation Focus (PS-AF)" level
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = "P2E_June_2025, P2E_July_2025, P2E_Sept_2025"
$ws.Range("H4").Value = "Affiliation Focus (AF), Merit Focus (MF), Personal Safety (PS), Search vs. Stay (SS), Affiliation Focus / Merit Focus (AF-MF), Personal Safety / Affiliation Focus (PS-AF)"

# The longer text in H4 reflows within the (unchanged) column width, so the
# author's workbook ended up with slightly shorter autofit row heights for
# the header rows once resaved. Match those explicit row heights.
$ws.Rows.Item(1).RowHeight = 31.5
$ws.Rows.Item(2).RowHeight = 44.25
$ws.Rows.Item(3).RowHeight = 171.75
$ws.Rows.Item(4).RowHeight = 146.25
